$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new blank row above row 19 on the "Sale 22-23" sheet, shifting the
# existing rows 19-24 down to 20-25 (formulas auto-adjust on insert).
$ws2.Rows.Item(19).Insert()

# The newly inserted row should pick up the formatting of the row above it
# (row 18) for columns E:F, matching the "box" the row sits inside, rather
# than the blank spacer-row formatting used below.
$ws2.Range("E18:F18").Copy()
$ws2.Range("E19:F19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the saved selections on each sheet.
$ws1.Select()
$ws1.Range("E44").Select()

$ws2.Select()
$ws2.Range("E18").Select()
